$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$pdSheet = $wb.Worksheets.Item("PDiCECpDoC")

# Insert the new "Texas Notes" sheet between "About" and "PDiCECpDoC"
$txSheet = $wb.Worksheets.Add($pdSheet)
$txSheet.Name = "Texas Notes"

# Re-fetch sheets by name since the collection shifted after the insert
$aboutSheet = $wb.Worksheets.Item("About")
$pdSheet = $wb.Worksheets.Item("PDiCECpDoC")

# --- Populate Texas Notes sheet ---
$txSheet.Columns.Item(1).ColumnWidth = 16.75

$txSheet.Range("A1").Value = "A more recent study from DNVGL "
$txSheet.Range("A2").Value = "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html"
$txSheet.Hyperlinks.Add($txSheet.Range("A2"), "https://www.dnvgl.com/feature/carbon-capture-storage-ccs.html") | Out-Null
$txSheet.Range("A3").Value = "assumes a learning rate of 15-20% - closer to what we see in the renewables and storage industries."

$txSheet.Range("A5").Value = "This suggests that the learning rate may be higher than the 2013 report used by EI. "
$txSheet.Range("A6").Value = "So, we can take an average of some of these values just to be conservative"

$txSheet.Range("B7").Value = "average learning rate"

$txSheet.Range("A8").Value = "2018 DNVGL"
$txSheet.Range("A8").HorizontalAlignment = -4131
$txSheet.Range("B8").Formula = "=AVERAGE(0.15, 0.2)"

$txSheet.Range("A9").Value = "2013 CRS report"
$txSheet.Range("A9").HorizontalAlignment = -4131
$txSheet.Range("B9").Value = 0.13

$txSheet.Range("A10").Value = "average"
$txSheet.Range("A10").HorizontalAlignment = -4152
$txSheet.Range("B10").Formula = "=AVERAGE(B8:B9)"
$txSheet.Range("B10").Interior.Color = 65535

$txSheet.Range("G13").Select()

# --- Update PDiCECpDoC sheet ---
$pdSheet.Range("B2").Formula = "='Texas Notes'!B10"
$pdSheet.Range("C7").Select()

# --- Update About sheet selection ---
$aboutSheet.Range("E24").Select()

$pdSheet.Activate()
